$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 178.53334
$ws.Range("I28").Value = 191.27777
$ws.Range("J28").Value = 159.41667
$ws.Range("K28").Value = 191.27777
$ws.Range("L28").Value = 159.41667
$ws.Range("M28").Value = 293.72223
$ws.Range("N28").Value = -1129.41667

$ws.Range("H137").Value = 1182.7709
$ws.Range("I137").Value = 1058.7949
$ws.Range("K137").Value = 3176.384700000001
$ws.Range("M137").Value = -626.3847000000005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 15000500
$ws.Range("I8").Value = 15000500
$ws.Range("K8").Value = 15000500
$ws.Range("M8").Value = -15000356

$ws.Range("H32").Value = 1941.43
$ws.Range("I32").Value = 2026.2273
$ws.Range("J32").Value = 1319.5834
$ws.Range("K32").Value = 2026.2273
$ws.Range("L32").Value = 1319.5834
$ws.Range("M32").Value = -1739.2273
$ws.Range("N32").Value = -1893.5834

$ws.Range("H61").Value = 949.375
$ws.Range("I61").Value = 885.9429
$ws.Range("J61").Value = 1120.1538
$ws.Range("K61").Value = 885.9429
$ws.Range("L61").Value = 1120.1538
$ws.Range("M61").Value = -673.9429
$ws.Range("N61").Value = -1544.1538

$ws.Range("H86").Value = 508900
$ws.Range("J86").Value = 1000000
$ws.Range("L86").Value = 1000000
$ws.Range("N86").Value = -1002372

$ws.Range("H89").Value = 508900
$ws.Range("J89").Value = 1000000
$ws.Range("L89").Value = 3000000
$ws.Range("N89").Value = -3011856

$ws.Range("H97").Value = 750.2619
$ws.Range("I97").Value = 622.84375
$ws.Range("J97").Value = 1158
$ws.Range("K97").Value = 622.84375
$ws.Range("L97").Value = 1158
$ws.Range("M97").Value = -126.84375
$ws.Range("N97").Value = -2150

$ws.Range("H136").Value = 949.375
$ws.Range("I136").Value = 885.9429
$ws.Range("J136").Value = 1120.1538
$ws.Range("K136").Value = 2657.8287
$ws.Range("L136").Value = 3360.4614
$ws.Range("M136").Value = -107.8287
$ws.Range("N136").Value = -8460.4614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 1153.75
$ws.Range("I29").Value = 871.6667
$ws.Range("K29").Value = 871.6667
$ws.Range("M29").Value = -582.6667

$ws.Range("H86").Value = 1012797.2
$ws.Range("I86").Value = 2173.9
$ws.Range("J86").Value = 1790199.8
$ws.Range("K86").Value = 2173.9
$ws.Range("L86").Value = 1790199.8
$ws.Range("M86").Value = -1050.9
$ws.Range("N86").Value = -1792445.8

$ws.Range("H89").Value = 1012797.2
$ws.Range("I89").Value = 2173.9
$ws.Range("J89").Value = 1790199.8
$ws.Range("K89").Value = 10869.5
$ws.Range("L89").Value = 8950999
$ws.Range("M89").Value = -5253.5
$ws.Range("N89").Value = -8962231

$ws.Range("H134").Value = 2927909.8
$ws.Range("I134").Value = 1243.8422
$ws.Range("J134").Value = 5854575.5
$ws.Range("K134").Value = 3731.5266
$ws.Range("L134").Value = 17563726.5
$ws.Range("M134").Value = -1196.5266
$ws.Range("N134").Value = -17568796.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1382.9756
$ws.Range("I31").Value = 1040.55
$ws.Range("J31").Value = 1709.0952
$ws.Range("K31").Value = 1040.55
$ws.Range("L31").Value = 1709.0952
$ws.Range("M31").Value = -745.55
$ws.Range("N31").Value = -2299.0952

$ws.Range("H34").Value = 1382.9756
$ws.Range("I34").Value = 1040.55
$ws.Range("J34").Value = 1709.0952
$ws.Range("K34").Value = 1040.55
$ws.Range("L34").Value = 1709.0952
$ws.Range("M34").Value = -838.55
$ws.Range("N34").Value = -2113.0952

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 36022.223
$ws.Range("J36").Value = 80725
$ws.Range("L36").Value = 242175
$ws.Range("N36").Value = -242513

$ws.Range("H80").Value = 1286.5
$ws.Range("I80").Value = 500
$ws.Range("J80").Value = 1548.6666
$ws.Range("K80").Value = 1500
$ws.Range("L80").Value = 4645.9998
$ws.Range("M80").Value = -564
$ws.Range("N80").Value = -6517.9998

$ws.Range("H83").Value = 1286.5
$ws.Range("I83").Value = 500
$ws.Range("J83").Value = 1548.6666
$ws.Range("K83").Value = 4500
$ws.Range("L83").Value = 13937.9994
$ws.Range("M83").Value = 180
$ws.Range("N83").Value = -23297.9994

$ws.Range("H131").Value = 818.8659699999999
$ws.Range("J131").Value = 866.52325
$ws.Range("L131").Value = 2599.56975
$ws.Range("N131").Value = -12679.56975

$ws.Range("H137").Value = 41667716
$ws.Range("J137").Value = 55556764
$ws.Range("L137").Value = 166670292
$ws.Range("N137").Value = -166680492

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 14051
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws.Range("H80").Value = 9094500
$ws.Range("I80").Value = 4438.125
$ws.Range("J80").Value = 33334666
$ws.Range("K80").Value = 4438.125
$ws.Range("L80").Value = 33334666
$ws.Range("M80").Value = -3440.125
$ws.Range("N80").Value = -33336662

$ws.Range("H83").Value = 9094500
$ws.Range("I83").Value = 4438.125
$ws.Range("J83").Value = 33334666
$ws.Range("K83").Value = 22190.625
$ws.Range("L83").Value = 166673330
$ws.Range("M83").Value = -17198.625
$ws.Range("N83").Value = -166683314

$ws.Range("H122").Value = 21744816
$ws.Range("I122").Value = 27784476
$ws.Range("J122").Value = 2040
$ws.Range("K122").Value = 83353428
$ws.Range("L122").Value = 6120
$ws.Range("M122").Value = -83350978
$ws.Range("N122").Value = -11020

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1067.8667
$ws.Range("I93").Value = 901.2
$ws.Range("K93").Value = 901.2
$ws.Range("M93").Value = 346.8

$ws.Range("H122").Value = 16642.857
$ws.Range("I122").Value = 21860
$ws.Range("J122").Value = 3600
$ws.Range("K122").Value = 65580
$ws.Range("L122").Value = 10800
$ws.Range("M122").Value = -63130
$ws.Range("N122").Value = -15700

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 21346.654
$ws.Range("J122").Value = 2396.9333
$ws.Range("L122").Value = 7190.7999
$ws.Range("N122").Value = -12090.7999

$ws.Range("H132").Value = 22395.926
$ws.Range("I132").Value = 24722.273
$ws.Range("J132").Value = 12160
$ws.Range("K132").Value = 74166.819
$ws.Range("L132").Value = 36480
$ws.Range("M132").Value = -71636.819
$ws.Range("N132").Value = -41540
